$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$albuns = $wb.Worksheets.Add($null, $ws1)
$albuns.Name = "Albuns"
$albuns.Cells.Item(1,1).Value = "Album"
$nm = $albuns.Names.Add("_xlnm._FilterDatabase", "=Albuns!`$A`$1:`$B`$11")
try { $nm.Hidden = $true } catch { Write-Output "Hidden failed: $_" }
try { Write-Output $nm.Visible } catch { Write-Output "get Visible failed: $_" }
Write-Output "done"
